# Companies Changes 1 - 29th June 2023
#
# The "Users" sheet holds the single test-user name used by the
# CapIQCompanies_CapIQCompaniesDetailPage_AddSalesforceCompany test.
# Swap the outgoing user for the new one and leave it as the
# active sheet/cell, matching how it was left after editing in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Activate()
$ws.Range("A2").Value = "Drew Koecher"
$ws.Range("A2").Select()
